$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2023-09-24 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-09-25 Monday", 2)

# Update the division problems in the table, cell by cell (row, column)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "97÷5="
$t.Cell(1, 2).Range.Text = "70÷7="
$t.Cell(1, 3).Range.Text = "63÷9="
$t.Cell(1, 4).Range.Text = "74÷6="
$t.Cell(1, 5).Range.Text = "67÷2="

$t.Cell(5, 1).Range.Text = "45÷7="
$t.Cell(5, 2).Range.Text = "13÷3="
$t.Cell(5, 3).Range.Text = "28÷3="
$t.Cell(5, 4).Range.Text = "76÷2="
$t.Cell(5, 5).Range.Text = "93÷5="

$t.Cell(9, 1).Range.Text = "31÷7="
$t.Cell(9, 2).Range.Text = "42÷3="
$t.Cell(9, 3).Range.Text = "79÷2="
$t.Cell(9, 4).Range.Text = "69÷8="
$t.Cell(9, 5).Range.Text = "68÷2="

$t.Cell(13, 1).Range.Text = "73÷9="
$t.Cell(13, 2).Range.Text = "31÷6="
$t.Cell(13, 3).Range.Text = "33÷3="
$t.Cell(13, 4).Range.Text = "69÷4="
$t.Cell(13, 5).Range.Text = "42÷3="

$t.Cell(17, 1).Range.Text = "91÷2="
$t.Cell(17, 2).Range.Text = "62÷6="
$t.Cell(17, 3).Range.Text = "94÷8="
$t.Cell(17, 4).Range.Text = "64÷9="
$t.Cell(17, 5).Range.Text = "18÷5="
